$wb = $excel.ActiveWorkbook

# Sheet ALC, row 13 (@@ -1278,19 +1278,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 5000
$ws.Range("J13").Value = 5000
$ws.Range("L13").Value = -5338

# Sheet ALC, row 16 (@@ -1419,19 +1422,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 2449.6667
$ws.Range("I16").Value = 2449.6667
$ws.Range("K16").Value = 2449.6667
$ws.Range("L16").Value = -2219.6667

# Sheet ALC, row 21 (@@ -1664,19 +1670,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 5569.375
$ws.Range("I21").Value = 3000
$ws.Range("J21").Value = 5936.4287
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = -6872.4287

# Sheet ALC, row 23 (@@ -1759,19 +1771,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 5569.375
$ws.Range("I23").Value = 3000
$ws.Range("J23").Value = 5936.4287
$ws.Range("K23").Value = 3000
$ws.Range("L23").Value = -6404.4287

# Sheet ALC, row 28 (@@ -1995,25 +2013,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 475.0625
$ws.Range("I28").Value = 515.75
$ws.Range("J28").Value = 353
$ws.Range("K28").Value = 515.75
$ws.Range("L28").Value = 353
$ws.Range("M28").Value = -30.75
$ws.Range("N28").Value = -1323

# Sheet ALC, row 64 (@@ -3813,25 +3831,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4057.647
$ws.Range("I64").Value = 3768.3684
$ws.Range("J64").Value = 4424.067
$ws.Range("K64").Value = 3768.3684
$ws.Range("L64").Value = 4424.067
$ws.Range("M64").Value = -3520.3684
$ws.Range("N64").Value = -4920.067

# Sheet ALC, row 67 (@@ -3969,25 +3987,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4057.647
$ws.Range("I67").Value = 3768.3684
$ws.Range("J67").Value = 4424.067
$ws.Range("K67").Value = 3768.3684
$ws.Range("L67").Value = 4424.067
$ws.Range("M67").Value = -2910.3684
$ws.Range("N67").Value = -6140.067

# Sheet ALC, row 113 (@@ -6301,25 +6319,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3624.2334
$ws.Range("I113").Value = 3759.8333
$ws.Range("J113").Value = 3590.3333
$ws.Range("K113").Value = 3759.8333
$ws.Range("L113").Value = 3590.3333
$ws.Range("M113").Value = -505.8332999999998
$ws.Range("N113").Value = -10098.3333

# Sheet ALC, row 116 (@@ -6451,22 +6469,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 39301.793
$ws.Range("I116").Value = 64732.824
$ws.Range("K116").Value = 64732.824
$ws.Range("M116").Value = -61290.824

# Sheet ALC, row 132 (@@ -7259,25 +7277,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3015.3076
$ws.Range("I132").Value = 1933.2424
$ws.Range("J132").Value = 8966.666999999999
$ws.Range("K132").Value = 5799.7272
$ws.Range("L132").Value = 26900.001
$ws.Range("M132").Value = -3269.7272
$ws.Range("N132").Value = -31960.001

# Sheet ALC, row 137 (@@ -7507,25 +7525,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2425.5
$ws.Range("I137").Value = 2364.4517
$ws.Range("J137").Value = 2597.5454
$ws.Range("K137").Value = 7093.355100000001
$ws.Range("L137").Value = 7792.6362
$ws.Range("M137").Value = -4543.355100000001
$ws.Range("N137").Value = -12892.6362

# Sheet ARM, row 45 (@@ -9968,25 +9986,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1221.3182
$ws.Range("I45").Value = 927.5
$ws.Range("J45").Value = 1573.9
$ws.Range("K45").Value = 927.5
$ws.Range("L45").Value = 1573.9
$ws.Range("M45").Value = -550.5
$ws.Range("N45").Value = -2327.9

# Sheet ARM, row 61 (@@ -10737,25 +10755,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 639301.5600000001
$ws.Range("I61").Value = 581245
$ws.Range("J61").Value = 718092.6
$ws.Range("K61").Value = 581245
$ws.Range("L61").Value = 718092.6
$ws.Range("M61").Value = -581033
$ws.Range("N61").Value = -718516.6

# Sheet ARM, row 74 (@@ -11371,25 +11389,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 181862.36
$ws.Range("I74").Value = 209352.05
$ws.Range("J74").Value = 80362
$ws.Range("K74").Value = 209352.05
$ws.Range("L74").Value = 80362
$ws.Range("M74").Value = -208478.05
$ws.Range("N74").Value = -82110

# Sheet ARM, row 77 (@@ -11521,25 +11539,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 181862.36
$ws.Range("I77").Value = 209352.05
$ws.Range("J77").Value = 80362
$ws.Range("K77").Value = 1046760.25
$ws.Range("L77").Value = 401810
$ws.Range("M77").Value = -1042392.25
$ws.Range("N77").Value = -410546

# Sheet ARM, row 88 (@@ -12072,25 +12090,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2244.074
$ws.Range("I88").Value = 2064.1177
$ws.Range("J88").Value = 2550
$ws.Range("K88").Value = 2064.1177
$ws.Range("L88").Value = 2550
$ws.Range("M88").Value = -1658.1177
$ws.Range("N88").Value = -3362

# Sheet ARM, row 91 (@@ -12219,25 +12237,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2244.074
$ws.Range("I91").Value = 2064.1177
$ws.Range("J91").Value = 2550
$ws.Range("K91").Value = 2064.1177
$ws.Range("L91").Value = 2550
$ws.Range("M91").Value = -660.1176999999998
$ws.Range("N91").Value = -5358

# Sheet ARM, row 136 (@@ -14418,25 +14436,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 639301.5600000001
$ws.Range("I136").Value = 581245
$ws.Range("J136").Value = 718092.6
$ws.Range("K136").Value = 1743735
$ws.Range("L136").Value = 2154277.8
$ws.Range("M136").Value = -1741185
$ws.Range("N136").Value = -2159377.8

# Sheet BSM, row 105 (@@ -19895,25 +19913,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3847932
$ws.Range("I105").Value = 1751.6666
$ws.Range("J105").Value = 7144658
$ws.Range("K105").Value = 1751.6666
$ws.Range("L105").Value = 7144658
$ws.Range("M105").Value = -4.666600000000017
$ws.Range("N105").Value = -7148152

# Sheet CRP, row 31 (@@ -23217,25 +23235,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2164.4333
$ws.Range("I31").Value = 1442.2858
$ws.Range("J31").Value = 5381.273
$ws.Range("K31").Value = 1442.2858
$ws.Range("L31").Value = 5381.273
$ws.Range("M31").Value = -1147.2858
$ws.Range("N31").Value = -5971.273

# Sheet CRP, row 34 (@@ -23370,25 +23388,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2164.4333
$ws.Range("I34").Value = 1442.2858
$ws.Range("J34").Value = 5381.273
$ws.Range("K34").Value = 1442.2858
$ws.Range("L34").Value = 5381.273
$ws.Range("M34").Value = -1240.2858
$ws.Range("N34").Value = -5785.273

# Sheet CRP, row 58 (@@ -24540,22 +24558,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4357.763
$ws.Range("I58").Value = 5187.2915
$ws.Range("K58").Value = 5187.2915
$ws.Range("M58").Value = -4984.2915

# Sheet CRP, row 132 (@@ -28190,25 +28208,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1841.0682
$ws.Range("I132").Value = 1173.5938
$ws.Range("J132").Value = 3621
$ws.Range("K132").Value = 3520.7814
$ws.Range("L132").Value = 10863
$ws.Range("M132").Value = -990.7814000000003
$ws.Range("N132").Value = -15923

# Sheet CRP, row 136 (@@ -28395,22 +28413,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4357.763
$ws.Range("I136").Value = 5187.2915
$ws.Range("K136").Value = 15561.8745
$ws.Range("M136").Value = -13011.8745

# Sheet CUL, row 56 (@@ -31528,22 +31546,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4633.0835
$ws.Range("I56").Value = 4633.0835
$ws.Range("K56").Value = 4633.0835
$ws.Range("M56").Value = -4103.0835

# Sheet CUL, row 113 (@@ -34423,25 +34441,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 9091437
$ws.Range("I113").Value = 13889414
$ws.Range("J113").Value = 533.0526
$ws.Range("K113").Value = 41668242
$ws.Range("L113").Value = 1599.1578
$ws.Range("M113").Value = -41666072
$ws.Range("N113").Value = -5939.1578

# Sheet CUL, row 121 (@@ -34836,25 +34854,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 3655.1667
$ws.Range("I121").Value = 8482.5
$ws.Range("J121").Value = 2689.7
$ws.Range("K121").Value = 25447.5
$ws.Range("L121").Value = 8069.099999999999
$ws.Range("M121").Value = -24137.5
$ws.Range("N121").Value = -10689.1

# Sheet GSM, row 113 (@@ -41428,25 +41446,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1186.3529
$ws.Range("I113").Value = 919.0769
$ws.Range("J113").Value = 2055
$ws.Range("K113").Value = 919.0769
$ws.Range("L113").Value = 2055
$ws.Range("M113").Value = 1250.9231
$ws.Range("N113").Value = -6395

# Sheet GSM, row 122 (@@ -41869,22 +41887,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I122").Value = 1479.5
$ws.Range("J122").Value = 1219.2307
$ws.Range("K122").Value = 4438.5
$ws.Range("L122").Value = 3657.6921
$ws.Range("M122").Value = -1988.5
$ws.Range("N122").Value = -8557.6921

# Sheet GSM, row 126 (@@ -42062,25 +42080,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2989.3225
$ws.Range("I126").Value = 2752.5833
$ws.Range("J126").Value = 3138.842
$ws.Range("K126").Value = 8257.749899999999
$ws.Range("L126").Value = 9416.526
$ws.Range("M126").Value = -5787.749899999999
$ws.Range("N126").Value = -14356.526

# Sheet GSM, row 132 (@@ -42359,25 +42377,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4008.5518
$ws.Range("I132").Value = 3967.3044
$ws.Range("J132").Value = 4166.6665
$ws.Range("K132").Value = 11901.9132
$ws.Range("L132").Value = 12499.9995
$ws.Range("M132").Value = -9371.913199999999
$ws.Range("N132").Value = -17559.9995

# Sheet GSM, row 140 (@@ -42757,22 +42775,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 29278.334
$ws.Range("J140").Value = 29278.334
$ws.Range("L140").Value = 29278.334
$ws.Range("N140").Value = -39638.334

# Sheet LTW, row 122 (@@ -48898,25 +48916,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6731.8096
$ws.Range("I122").Value = 7547.9443
$ws.Range("J122").Value = 1835
$ws.Range("K122").Value = 22643.8329
$ws.Range("L122").Value = 5505
$ws.Range("M122").Value = -20193.8329
$ws.Range("N122").Value = -10405

# Sheet LTW, row 132 (@@ -49394,25 +49412,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6482.841
$ws.Range("I132").Value = 1777.2142
$ws.Range("J132").Value = 14717.6875
$ws.Range("K132").Value = 5331.642599999999
$ws.Range("L132").Value = 44153.0625
$ws.Range("M132").Value = -2801.642599999999
$ws.Range("N132").Value = -49213.0625

# Sheet WVR, row 131 (@@ -56362,22 +56380,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 39800
$ws.Range("J131").Value = 39800
$ws.Range("L131").Value = 39800
$ws.Range("N131").Value = -49880
